$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.389221
$ws.Range("H2").Value = 4.167663
$ws.Range("I2").Value = 0.2910270461264192
$ws.Range("J2").Value = 0.2910270461264192
$ws.Range("M2").Value = 2.598166333333333
$ws.Range("N2").Value = 7.794499
$ws.Range("O2").Value = 0.3466013321552429
$ws.Range("P2").Value = 0.3466013321552429
$ws.Range("Q2").Value = 3.609427231759667
$ws.Range("R2").Value = 32.484845085837
$ws.Range("S2").Value = 0.1008703618806222
$ws.Range("T2").Value = 0.1008703618806222
$ws.Range("G3").Value = 1.389221
$ws.Range("H3").Value = 4.167663
$ws.Range("I3").Value = 0.2910270461264192
$ws.Range("J3").Value = 0.2910270461264192
$ws.Range("M3").Value = 4.333403333333333
$ws.Range("O3").Value = 0.5780859172985858
$ws.Range("P3").Value = 0.5780859172985858
$ws.Range("Q3").Value = 6.020054912136667
$ws.Range("R3").Value = 54.18049420923
$ws.Range("S3").Value = 0.1682386369186889
$ws.Range("T3").Value = 0.1682386369186889
$ws.Range("G4").Value = 1.389221
$ws.Range("H4").Value = 4.167663
$ws.Range("I4").Value = 0.2910270461264192
$ws.Range("J4").Value = 0.2910270461264192
$ws.Range("M4").Value = 0.4692043333333333
$ws.Range("N4").Value = 1.407613
$ws.Range("O4").Value = 0.06259293136852516
$ws.Range("P4").Value = 0.06259293136852516
$ws.Range("Q4").Value = 0.6518285131576667
$ws.Range("R4").Value = 5.866456618419
$ws.Range("S4").Value = 0.01821623592457556
$ws.Range("T4").Value = 0.01821623592457556
$ws.Range("G5").Value = 1.389221
$ws.Range("H5").Value = 4.167663
$ws.Range("I5").Value = 0.2910270461264192
$ws.Range("J5").Value = 0.2910270461264192
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.09534933333333333
$ws.Range("N5").Value = 0.286048
$ws.Range("O5").Value = 0.01271981917764605
$ws.Range("P5").Value = 0.01271981917764604
$ws.Range("Q5").Value = 0.1324612962026667
$ws.Range("R5").Value = 1.192151665824
$ws.Range("S5").Value = 0.003701811402532507
$ws.Range("T5").Value = 0.003701811402532506
$ws.Range("I6").Value = 0.461328155686921
$ws.Range("J6").Value = 0.4613281556869209
$ws.Range("M6").Value = 2.598166333333333
$ws.Range("N6").Value = 7.794499
$ws.Range("O6").Value = 0.3466013321552429
$ws.Range("P6").Value = 0.3466013321552429
$ws.Range("Q6").Value = 5.721565847837112
$ws.Range("R6").Value = 51.49409263053401
$ws.Range("S6").Value = 0.1598969533218081
$ws.Range("T6").Value = 0.1598969533218081
$ws.Range("I7").Value = 0.461328155686921
$ws.Range("J7").Value = 0.4613281556869209
$ws.Range("M7").Value = 4.333403333333333
$ws.Range("O7").Value = 0.5780859172985858
$ws.Range("P7").Value = 0.5780859172985858
$ws.Range("Q7").Value = 9.542827261984446
$ws.Range("R7").Value = 85.88544535786001
$ws.Range("S7").Value = 0.2666873100559385
$ws.Range("T7").Value = 0.2666873100559385
$ws.Range("I8").Value = 0.461328155686921
$ws.Range("J8").Value = 0.4613281556869209
$ws.Range("M8").Value = 0.4692043333333333
$ws.Range("N8").Value = 1.407613
$ws.Range("O8").Value = 0.06259293136852516
$ws.Range("P8").Value = 0.06259293136852516
$ws.Range("Q8").Value = 1.033260825073111
$ws.Range("R8").Value = 9.299347425658002
$ws.Range("S8").Value = 0.02887588158727973
$ws.Range("T8").Value = 0.02887588158727973
$ws.Range("I9").Value = 0.461328155686921
$ws.Range("J9").Value = 0.4613281556869209
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.09534933333333333
$ws.Range("N9").Value = 0.286048
$ws.Range("O9").Value = 0.01271981917764605
$ws.Range("P9").Value = 0.01271981917764604
$ws.Range("Q9").Value = 0.2099740429297778
$ws.Range("R9").Value = 1.889766386368
$ws.Range("S9").Value = 0.005868010721894578
$ws.Range("T9").Value = 0.005868010721894577
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.1506176666666667
$ws.Range("H10").Value = 0.451853
$ws.Range("I10").Value = 0.03155280162368235
$ws.Range("J10").Value = 0.03155280162368235
$ws.Range("M10").Value = 2.598166333333333
$ws.Range("N10").Value = 7.794499
$ws.Range("O10").Value = 0.3466013321552429
$ws.Range("P10").Value = 0.3466013321552429
$ws.Range("Q10").Value = 0.3913297507385556
$ws.Range("R10").Value = 3.521967756647
$ws.Range("S10").Value = 0.01093624307599842
$ws.Range("T10").Value = 0.01093624307599842
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.1506176666666667
$ws.Range("H11").Value = 0.451853
$ws.Range("I11").Value = 0.03155280162368235
$ws.Range("J11").Value = 0.03155280162368235
$ws.Range("M11").Value = 4.333403333333333
$ws.Range("O11").Value = 0.5780859172985858
$ws.Range("P11").Value = 0.5780859172985858
$ws.Range("Q11").Value = 0.6526870987922222
$ws.Range("R11").Value = 5.874183889129999
$ws.Range("S11").Value = 0.01824023026996672
$ws.Range("T11").Value = 0.01824023026996672
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.1506176666666667
$ws.Range("H12").Value = 0.451853
$ws.Range("I12").Value = 0.03155280162368235
$ws.Range("J12").Value = 0.03155280162368235
$ws.Range("M12").Value = 0.4692043333333333
$ws.Range("N12").Value = 1.407613
$ws.Range("O12").Value = 0.06259293136852516
$ws.Range("P12").Value = 0.06259293136852516
$ws.Range("Q12").Value = 0.07067046187655555
$ws.Range("R12").Value = 0.636034156889
$ws.Range("S12").Value = 0.001974982346515839
$ws.Range("T12").Value = 0.001974982346515839
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.1506176666666667
$ws.Range("H13").Value = 0.451853
$ws.Range("I13").Value = 0.03155280162368235
$ws.Range("J13").Value = 0.03155280162368235
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.09534933333333333
$ws.Range("N13").Value = 0.286048
$ws.Range("O13").Value = 0.01271981917764605
$ws.Range("P13").Value = 0.01271981917764604
$ws.Range("Q13").Value = 0.01436129410488889
$ws.Range("R13").Value = 0.129251646944
$ws.Range("S13").Value = 0.0004013459312013761
$ws.Range("T13").Value = 0.000401345931201376
$ws.Range("G14").Value = 1.031517666666667
$ws.Range("H14").Value = 3.094553
$ws.Range("I14").Value = 0.2160919965629775
$ws.Range("J14").Value = 0.2160919965629775
$ws.Range("M14").Value = 2.598166333333333
$ws.Range("N14").Value = 7.794499
$ws.Range("O14").Value = 0.3466013321552429
$ws.Range("P14").Value = 0.3466013321552429
$ws.Range("Q14").Value = 2.680054473771889
$ws.Range("R14").Value = 24.120490263947
$ws.Range("S14").Value = 0.07489777387681419
$ws.Range("T14").Value = 0.07489777387681419
$ws.Range("G15").Value = 1.031517666666667
$ws.Range("H15").Value = 3.094553
$ws.Range("I15").Value = 0.2160919965629775
$ws.Range("J15").Value = 0.2160919965629775
$ws.Range("M15").Value = 4.333403333333333
$ws.Range("O15").Value = 0.5780859172985858
$ws.Range("P15").Value = 0.5780859172985858
$ws.Range("Q15").Value = 4.469982095125554
$ws.Range("R15").Value = 40.22983885612999
$ws.Range("S15").Value = 0.1249197400539917
$ws.Range("T15").Value = 0.1249197400539917
$ws.Range("G16").Value = 1.031517666666667
$ws.Range("H16").Value = 3.094553
$ws.Range("I16").Value = 0.2160919965629775
$ws.Range("J16").Value = 0.2160919965629775
$ws.Range("M16").Value = 0.4692043333333333
$ws.Range("N16").Value = 1.407613
$ws.Range("O16").Value = 0.06259293136852516
$ws.Range("P16").Value = 0.06259293136852516
$ws.Range("Q16").Value = 0.4839925591098889
$ws.Range("R16").Value = 4.355933031989
$ws.Range("S16").Value = 0.01352583151015403
$ws.Range("T16").Value = 0.01352583151015403
$ws.Range("G17").Value = 1.031517666666667
$ws.Range("H17").Value = 3.094553
$ws.Range("I17").Value = 0.2160919965629775
$ws.Range("J17").Value = 0.2160919965629775
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.09534933333333333
$ws.Range("N17").Value = 0.286048
$ws.Range("O17").Value = 0.01271981917764605
$ws.Range("P17").Value = 0.01271981917764604
$ws.Range("Q17").Value = 0.0983545218382222
$ws.Range("R17").Value = 0.8851906965439998
$ws.Range("S17").Value = 0.002748651122017585
$ws.Range("T17").Value = 0.002748651122017585
